$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.661.19"
$ws.Range("E2").Value = "  +5.72%  "

$ws.Range("D3").Value = "3.051.96"
$ws.Range("E3").Value = "  +5.26%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "555.86"
$ws.Range("E5").Value = "  +5.05%  "

$ws.Range("D6").Value = "142.01"
$ws.Range("E6").Value = "  +8.58%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "3.045.66"
$ws.Range("E8").Value = "  +5.50%  "

$ws.Range("E9").Value = "  +6.62%  "

$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  +9.77%  "

$ws.Range("E11").Value = "  -3.06%  "

$ws.Range("D12").Value = "0.478"
$ws.Range("E12").Value = "  +11.57%  "

$ws.Range("E13").Value = "  +9.09%  "

$ws.Range("D14").Value = "34.91"
$ws.Range("E14").Value = "  +7.82%  "

$ws.Range("D15").Value = "3.548.44"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("D16").Value = "63.674.33"
$ws.Range("E16").Value = "  +5.65%  "

$ws.Range("E17").Value = "  +4.29%  "

$ws.Range("D18").Value = "3.050.48"
$ws.Range("E18").Value = "  +4.62%  "

$ws.Range("E19").Value = "  +5.69%  "

$ws.Range("D20").Value = "477.82"
$ws.Range("E20").Value = "  +5.72%  "

$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  +8.71%  "

$ws.Range("E22").Value = "  +7.90%  "

$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +10.62%  "

$ws.Range("D24").Value = "14.25"
$ws.Range("E24").Value = "  +19.30%  "

$ws.Range("D25").Value = "81.29"
$ws.Range("E25").Value = "  +5.48%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  +6.44%  "

$ws.Range("D28").Value = "7.93"
$ws.Range("E28").Value = "  +9.44%  "

$ws.Range("D29").Value = "2.03"
$ws.Range("E29").Value = "  +6.91%  "

$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").Value = "26.17"
$ws.Range("E31").Value = "  +7.00%  "

$ws.Range("E32").Value = "  +3.96%  "

$ws.Range("D33").Value = "2.43"
$ws.Range("E33").Value = "  +8.43%  "

$ws.Range("D34").Value = "5.63"
$ws.Range("E34").Value = "  +5.91%  "

$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  +10.84%  "

$ws.Range("D36").Value = "54.85"
$ws.Range("E36").Value = "  +3.60%  "

$ws.Range("E37").Value = "  +9.66%  "

$ws.Range("D38").Value = "442.93"
$ws.Range("E38").Value = "  +2.37%  "

$ws.Range("E39").Value = "  +4.87%  "

$ws.Range("E40").Value = "  +23.58%  "

$ws.Range("D41").Value = "2.957.26"
$ws.Range("E41").Value = "  +3.48%  "

$ws.Range("E42").Value = "  +6.80%  "

$ws.Range("D43").Value = "0.112"
$ws.Range("E43").Value = "  +1.24%  "

$ws.Range("D44").Value = "27.80"
$ws.Range("E44").Value = "  +9.64%  "

$ws.Range("E45").Value = "  +9.78%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +14.29%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "0.112"
$ws.Range("E48").Value = "  +7.15%  "

$ws.Range("D49").Value = "0.0₃0514"
$ws.Range("E49").Value = "  +10.28%  "

$ws.Range("D50").Value = "116.79"
$ws.Range("E50").Value = "  +4.83%  "

$ws.Range("E51").Value = "  +8.46%  "
